$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a daily price log with the newest entry in row 2 and older
# entries below. A new day's row is being inserted at the top (row 2),
# pushing every existing row (2..179) down by one (to 3..180).
#
# Do this as a single bulk copy: copy A2:F179 down to A3:F180. This moves
# both values and styles (so the new row 180 correctly inherits the s="3"
# / s="4" styles used throughout, instead of picking up a blank default
# style). Row 2 itself is left untouched by this copy (the source range
# starts at row 2, the destination at row 3), so the old row-2 values
# (price/description/etc.) remain there, ready to become the new row 2.
$ws.Range("A2:F179").Copy($ws.Range("A3:F180"))

# Row 2 keeps all of the old row 2's data (price unchanged from the
# previous entry) except for the date, which advances to the new day.
# Force text (leading apostrophe) so Excel doesn't reinterpret the
# dd-mm-yyyy-looking string as a real date value.
$ws.Cells.Item(2, 1).Value = "'07-12-2025"

# Note: Range.Copy does not relocate Hyperlink objects, so the existing
# Hyperlinks collection (F2..F118 -> rId1..rId117) is left exactly as it
# was, still anchored to the same rows/ids -- matching the source data,
# which keeps those hyperlink relationships untouched even though the
# displayed text in those cells has shifted. Only one new hyperlink needs
# to be added, for the newly extended F119, reusing the same target that
# F118 already pointed to (the boundary simply moves one row down).
$ws.Hyperlinks.Add($ws.Range("F119"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")
